# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback
# DateTime" (H2) values for the 96bed73a-... file row on both the zh-cn
# and de-de language report sheets, reflecting the newly generated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 02:53:29"
$wsZhCn.Range("H2").Value = "2016-03-24 02:53:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 02:53:34"
$wsDeDe.Range("H2").Value = "2016-03-24 02:54:01"
